# Apply updated route/day/temperature/mean15d/half/best/date table
# (Build explanation text; Update progress tracker)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Louisiana Plus"
$ws.Cells.Item(2, 2).Value = 97
$ws.Cells.Item(2, 3).Value = 0.74
$ws.Cells.Item(2, 4).Value = 0.72
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = $True
$ws.Cells.Item(2, 7).Value = "04-07"

$ws.Cells.Item(3, 1).Value = "Louisiana Plus"
$ws.Cells.Item(3, 2).Value = 337
$ws.Cells.Item(3, 3).Value = 0.68
$ws.Cells.Item(3, 4).Value = 0.67
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = $True
$ws.Cells.Item(3, 7).Value = "12-03"

$ws.Cells.Item(4, 1).Value = "Near DC"
$ws.Cells.Item(4, 2).Value = 137
$ws.Cells.Item(4, 3).Value = 0.71
$ws.Cells.Item(4, 4).Value = 0.74
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = $True
$ws.Cells.Item(4, 7).Value = "05-17"

$ws.Cells.Item(5, 1).Value = "Near DC"
$ws.Cells.Item(5, 2).Value = 297
$ws.Cells.Item(5, 3).Value = 0.77
$ws.Cells.Item(5, 4).Value = 0.79
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = $True
$ws.Cells.Item(5, 7).Value = "10-24"

$ws.Cells.Item(6, 1).Value = "Oregon Plus"
$ws.Cells.Item(6, 2).Value = 177
$ws.Cells.Item(6, 3).Value = 0.74
$ws.Cells.Item(6, 4).Value = 0.77
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = $True
$ws.Cells.Item(6, 7).Value = "06-26"

$ws.Cells.Item(7, 1).Value = "Oregon Plus"
$ws.Cells.Item(7, 2).Value = 272
$ws.Cells.Item(7, 3).Value = 0.82
$ws.Cells.Item(7, 4).Value = 0.81
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = $True
$ws.Cells.Item(7, 7).Value = "09-29"

$ws.Cells.Item(8, 1).Value = "Northeast Plus"
$ws.Cells.Item(8, 2).Value = 157
$ws.Cells.Item(8, 3).Value = 0.74
$ws.Cells.Item(8, 4).Value = 0.72
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = $True
$ws.Cells.Item(8, 7).Value = "06-06"

$ws.Cells.Item(9, 1).Value = "Northeast Plus"
$ws.Cells.Item(9, 2).Value = 287
$ws.Cells.Item(9, 3).Value = 0.78
$ws.Cells.Item(9, 4).Value = 0.77
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = $True
$ws.Cells.Item(9, 7).Value = "10-14"

$ws.Cells.Item(10, 1).Value = "Minnesota Plus"
$ws.Cells.Item(10, 2).Value = 177
$ws.Cells.Item(10, 3).Value = 0.8
$ws.Cells.Item(10, 4).Value = 0.71
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = $True
$ws.Cells.Item(10, 7).Value = "06-26"

$ws.Cells.Item(11, 1).Value = "Minnesota Plus"
$ws.Cells.Item(11, 2).Value = 252
$ws.Cells.Item(11, 3).Value = 0.77
$ws.Cells.Item(11, 4).Value = 0.77
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = $True
$ws.Cells.Item(11, 7).Value = "09-09"

$ws.Cells.Item(12, 1).Value = "Hawaii State"
$ws.Cells.Item(12, 2).Value = 52
$ws.Cells.Item(12, 3).Value = 0.76
$ws.Cells.Item(12, 4).Value = 0.79
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = $True
$ws.Cells.Item(12, 7).Value = "02-21"

$ws.Cells.Item(13, 1).Value = "Hawaii State"
$ws.Cells.Item(13, 2).Value = 365
$ws.Cells.Item(13, 3).Value = 0.64
$ws.Cells.Item(13, 4).Value = 0.63
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = $True
$ws.Cells.Item(13, 7).Value = "12-31"

$ws.Cells.Item(14, 1).Value = "Alaska State"
$ws.Cells.Item(14, 2).Value = 200
$ws.Cells.Item(14, 3).Value = 0.85
$ws.Cells.Item(14, 4).Value = 0.88
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = $True
$ws.Cells.Item(14, 7).Value = "07-19"

$ws.Cells.Item(15, 1).Value = "Alaska State"
$ws.Cells.Item(15, 2).Value = 246
$ws.Cells.Item(15, 3).Value = 0.87
$ws.Cells.Item(15, 4).Value = 0.87
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = $True
$ws.Cells.Item(15, 7).Value = "09-03"

$ws.Cells.Item(16, 1).Value = "Wyoming Plus"
$ws.Cells.Item(16, 2).Value = 157
$ws.Cells.Item(16, 3).Value = 0.66
$ws.Cells.Item(16, 4).Value = 0.65
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = $True
$ws.Cells.Item(16, 7).Value = "06-06"

$ws.Cells.Item(17, 1).Value = "Wyoming Plus"
$ws.Cells.Item(17, 2).Value = 272
$ws.Cells.Item(17, 3).Value = 0.61
$ws.Cells.Item(17, 4).Value = 0.62
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = $True
$ws.Cells.Item(17, 7).Value = "09-29"

$ws.Cells.Item(18, 1).Value = "Florida State"
$ws.Cells.Item(18, 2).Value = 37
$ws.Cells.Item(18, 3).Value = 0.71
$ws.Cells.Item(18, 4).Value = 0.7
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = $True
$ws.Cells.Item(18, 7).Value = "02-06"

$ws.Cells.Item(19, 1).Value = "Florida State"
$ws.Cells.Item(19, 2).Value = 362
$ws.Cells.Item(19, 3).Value = 0.69
$ws.Cells.Item(19, 4).Value = 0.67
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = $True
$ws.Cells.Item(19, 7).Value = "12-28"

$ws.Cells.Item(20, 1).Value = "Missouri Plus"
$ws.Cells.Item(20, 2).Value = 147
$ws.Cells.Item(20, 3).Value = 0.72
$ws.Cells.Item(20, 4).Value = 0.72
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = $True
$ws.Cells.Item(20, 7).Value = "05-27"

$ws.Cells.Item(21, 1).Value = "Missouri Plus"
$ws.Cells.Item(21, 2).Value = 282
$ws.Cells.Item(21, 3).Value = 0.74
$ws.Cells.Item(21, 4).Value = 0.69
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 6).Value = $True
$ws.Cells.Item(21, 7).Value = "10-09"

$ws.Cells.Item(22, 1).Value = "Puerto Rico"
$ws.Cells.Item(22, 2).Value = 67
$ws.Cells.Item(22, 3).Value = 0.1
$ws.Cells.Item(22, 4).Value = 0.08
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = $True
$ws.Cells.Item(22, 7).Value = "03-08"

$ws.Cells.Item(23, 1).Value = "Puerto Rico"
$ws.Cells.Item(23, 2).Value = 365
$ws.Cells.Item(23, 3).Value = 0.05
$ws.Cells.Item(23, 4).Value = 0.02
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = $True
$ws.Cells.Item(23, 7).Value = "12-31"

$ws.Cells.Item(24, 1).Value = "East Canada"
$ws.Cells.Item(24, 2).Value = 182
$ws.Cells.Item(24, 3).Value = 0.77
$ws.Cells.Item(24, 4).Value = 0.77
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = $True
$ws.Cells.Item(24, 7).Value = "07-01"

$ws.Cells.Item(25, 1).Value = "East Canada"
$ws.Cells.Item(25, 2).Value = 257
$ws.Cells.Item(25, 3).Value = 0.78
$ws.Cells.Item(25, 4).Value = 0.81
$ws.Cells.Item(25, 5).Value = 1
$ws.Cells.Item(25, 6).Value = $True
$ws.Cells.Item(25, 7).Value = "09-14"

$ws.Cells.Item(26, 1).Value = "California Plus"
$ws.Cells.Item(26, 2).Value = 142
$ws.Cells.Item(26, 3).Value = 0.72
$ws.Cells.Item(26, 4).Value = 0.73
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = $True
$ws.Cells.Item(26, 7).Value = "05-22"

$ws.Cells.Item(27, 1).Value = "California Plus"
$ws.Cells.Item(27, 2).Value = 312
$ws.Cells.Item(27, 3).Value = 0.71
$ws.Cells.Item(27, 4).Value = 0.73
$ws.Cells.Item(27, 5).Value = 1
$ws.Cells.Item(27, 6).Value = $True
$ws.Cells.Item(27, 7).Value = "11-08"

$ws.Cells.Item(28, 1).Value = "Georgia Plus"
$ws.Cells.Item(28, 2).Value = 117
$ws.Cells.Item(28, 3).Value = 0.78
$ws.Cells.Item(28, 4).Value = 0.75
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = $True
$ws.Cells.Item(28, 7).Value = "04-27"

$ws.Cells.Item(29, 1).Value = "Georgia Plus"
$ws.Cells.Item(29, 2).Value = 302
$ws.Cells.Item(29, 3).Value = 0.74
$ws.Cells.Item(29, 4).Value = 0.75
$ws.Cells.Item(29, 5).Value = 1
$ws.Cells.Item(29, 6).Value = $True
$ws.Cells.Item(29, 7).Value = "10-29"

$ws.Cells.Item(30, 1).Value = "Indiana Plus"
$ws.Cells.Item(30, 2).Value = 157
$ws.Cells.Item(30, 3).Value = 0.76
$ws.Cells.Item(30, 4).Value = 0.72
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = $True
$ws.Cells.Item(30, 7).Value = "06-06"

$ws.Cells.Item(31, 1).Value = "Indiana Plus"
$ws.Cells.Item(31, 2).Value = 262
$ws.Cells.Item(31, 3).Value = 0.73
$ws.Cells.Item(31, 4).Value = 0.75
$ws.Cells.Item(31, 5).Value = 1
$ws.Cells.Item(31, 6).Value = $True
$ws.Cells.Item(31, 7).Value = "09-19"

$ws.Cells.Item(32, 1).Value = "New Mexico Plus"
$ws.Cells.Item(32, 2).Value = 107
$ws.Cells.Item(32, 3).Value = 0.59
$ws.Cells.Item(32, 4).Value = 0.6
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = $True
$ws.Cells.Item(32, 7).Value = "04-17"

$ws.Cells.Item(33, 1).Value = "New Mexico Plus"
$ws.Cells.Item(33, 2).Value = 297
$ws.Cells.Item(33, 3).Value = 0.6
$ws.Cells.Item(33, 4).Value = 0.59
$ws.Cells.Item(33, 5).Value = 1
$ws.Cells.Item(33, 6).Value = $True
$ws.Cells.Item(33, 7).Value = "10-24"

